$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.796.19"
$ws.Range("E2").Value = "  -0.17%  "
$ws.Range("D3").Value = "2.077.70"
$ws.Range("E3").Value = "  -0.57%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "'232.99"
$ws.Range("E5").Value = "  -0.76%  "
$ws.Range("E6").Value = "  -0.21%  "
$ws.Range("D7").Value = "'58.54"
$ws.Range("E7").Value = "  -1.32%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  +0.31%  "
$ws.Range("D10").Value = "'0.0785"
$ws.Range("E10").Value = "  -1.04%  "
$ws.Range("E11").Value = "  +3.36%  "
$ws.Range("B12").Value = "Chainlink"
$ws.Range("C12").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D12").Value = "'14.86"
$ws.Range("E12").Value = "  +0.84%  "
$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").Value = "2.383.64"
$ws.Range("E13").Value = "  -0.53%  "
$ws.Range("E14").Value = "  -1.63%  "
$ws.Range("D15").Value = "'0.782"
$ws.Range("E15").Value = "  +1.16%  "
$ws.Range("D16").Value = "'5.35"
$ws.Range("E16").Value = "  +0.68%  "
$ws.Range("D17").Value = "2.098.92"
$ws.Range("E17").Value = "  +0.50%  "
$ws.Range("D18").Value = "37.687.25"
$ws.Range("D19").Value = "'6.13"
$ws.Range("E19").Value = "  -2.02%  "
$ws.Range("D20").Value = "'71.61"
$ws.Range("E20").Value = "  -0.41%  "
$ws.Range("D21").Value = "0.0₃0843"
$ws.Range("E21").Value = "  +1.41%  "
$ws.Range("D22").Value = "'229.20"
$ws.Range("E22").Value = "  +0.07%  "
$ws.Range("E24").Value = "  -1.09%  "
$ws.Range("E25").Value = "  -0.07%  "
$ws.Range("D26").Value = "'9.70"
$ws.Range("E26").Value = "  +6.76%  "
$ws.Range("D27").Value = "'171.77"
$ws.Range("E27").Value = "  +0.55%  "
$ws.Range("E28").Value = "  -1.09%  "
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").Value = "'19.43"
$ws.Range("E29").Value = "  -0.85%  "
$ws.Range("B30").Value = "ImmutableX"
$ws.Range("C30").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D30").Value = "'1.40"
$ws.Range("E30").Value = "  -2.14%  "
$ws.Range("E31").Value = "  +1.16%  "
$ws.Range("E32").Value = "  +0.53%  "
$ws.Range("E33").Value = "  -0.01%  "
$ws.Range("E34").Value = "  -0.99%  "
$ws.Range("E35").Value = "  -2.55%  "
$ws.Range("E36").Value = "  -0.19%  "
$ws.Range("D37").Value = "'3.40"
$ws.Range("E37").Value = "  -3.20%  "
$ws.Range("D38").Value = "'0.999"
$ws.Range("E38").Value = "  -0.08%  "
$ws.Range("D39").Value = "'5.43"
$ws.Range("E39").Value = "  -0.78%  "
$ws.Range("D40").Value = "'0.0231"
$ws.Range("E40").Value = "  +7.15%  "
$ws.Range("D41").Value = "'101.74"
$ws.Range("E41").Value = "  +1.99%  "
$ws.Range("D42").Value = "'0.0976"
$ws.Range("E42").Value = "  -1.22%  "
$ws.Range("E43").Value = "  -0.53%  "
$ws.Range("E44").Value = "  +4.77%  "
$ws.Range("D45").Value = "1.451.73"
$ws.Range("E45").Value = "  -0.81%  "
$ws.Range("E46").Value = "  -2.05%  "
$ws.Range("E47").Value = "  -1.67%  "
$ws.Range("E48").Value = "  -5.83%  "
$ws.Range("D49").Value = "'7.39"
$ws.Range("E49").Value = "  -1.42%  "
$ws.Range("E50").Value = "  -1.48%  "
$ws.Range("D51").Value = "2.269.77"
$ws.Range("E51").Value = "  -0.51%  "
